$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I: "MVP 2.0" task
$ws.Range("I1").Value = "MVP 2.0"
$ws.Range("I2").Value = "x"
$ws.Range("I3").Value = "x"
$ws.Range("I4").Value = "x"
$ws.Range("I5").Value = "x"
$ws.Range("I6").Value = "x"
$ws.Range("I8").Value = "Pütter/Hesse"
$ws.Range("I8").HorizontalAlignment = -4108

# Move the active selection to I9
$ws.Range("I9").Select()
